$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.196078431372549
$ws.Range("C2").Value = 0.5490196078431373
$ws.Range("J2").Value = 0.01143790849673203
$ws.Range("O2").Value = 0.001633986928104575
$ws.Range("P2").Value = 0.1470588235294118
$ws.Range("S2").Value = 0.09477124183006536
$ws.Range("B3").Value = 0.01652892561983471
$ws.Range("C3").Value = 0.06611570247933884
$ws.Range("J3").Value = 0.008264462809917356
$ws.Range("P3").Value = 0.721763085399449
$ws.Range("S3").Value = 0.1873278236914601
$ws.Range("J4").Value = 0.03896103896103896
$ws.Range("P4").Value = 0.6103896103896104
$ws.Range("S4").Value = 0.3506493506493507
$ws.Range("B6").Value = 0.07317073170731707
$ws.Range("D6").Value = 0.002217294900221729
$ws.Range("F6").Value = 0.06651884700665188
$ws.Range("J6").Value = 0.2483370288248337
$ws.Range("O6").Value = 0.01995565410199556
$ws.Range("Q6").Value = 0.1552106430155211
$ws.Range("R6").Value = 0.07760532150776053
$ws.Range("S6").Value = 0.3569844789356985
$ws.Range("B7").Value = 0.1111111111111111
$ws.Range("D7").Value = 0.0202020202020202
$ws.Range("E7").Value = 0.002525252525252525
$ws.Range("F7").Value = 0.05555555555555555
$ws.Range("J7").Value = 0.1237373737373737
$ws.Range("O7").Value = 0.01767676767676768
$ws.Range("Q7").Value = 0.1868686868686869
$ws.Range("R7").Value = 0.101010101010101
$ws.Range("S7").Value = 0.3813131313131313
$ws.Range("B8").Value = 0.09117361784675072
$ws.Range("D8").Value = 0.01745877788554801
$ws.Range("E8").Value = 0.0009699321047526673
$ws.Range("F8").Value = 0.05819592628516004
$ws.Range("J8").Value = 0.1134820562560621
$ws.Range("O8").Value = 0.01551891367604268
$ws.Range("Q8").Value = 0.1726479146459748
$ws.Range("R8").Value = 0.09990300678952474
$ws.Range("S8").Value = 0.4306498545101843
$ws.Range("B9").Value = 0.09429824561403509
$ws.Range("D9").Value = 0.008771929824561403
$ws.Range("E9").Value = 0.002192982456140351
$ws.Range("F9").Value = 0.06140350877192982
$ws.Range("J9").Value = 0.1074561403508772
$ws.Range("O9").Value = 0.01535087719298246
$ws.Range("Q9").Value = 0.1513157894736842
$ws.Range("R9").Value = 0.09429824561403509
$ws.Range("S9").Value = 0.4649122807017544
$ws.Range("B10").Value = 0.09939759036144578
$ws.Range("D10").Value = 0.01769578313253012
$ws.Range("E10").Value = 0.002259036144578313
$ws.Range("F10").Value = 0.06890060240963855
$ws.Range("J10").Value = 0.1118222891566265
$ws.Range("O10").Value = 0.01393072289156626
$ws.Range("Q10").Value = 0.1976656626506024
$ws.Range("R10").Value = 0.09186746987951808
$ws.Range("S10").Value = 0.396460843373494
$ws.Range("G11").Value = 0.1416234887737478
$ws.Range("J11").Value = 0.08117443868739206
$ws.Range("K11").Value = 0.1986183074265976
$ws.Range("L11").Value = 0.5630397236614854
$ws.Range("S11").Value = 0.0155440414507772
$ws.Range("G12").Value = 0.7492447129909365
$ws.Range("J12").Value = 0.1812688821752266
$ws.Range("K12").Value = 0.003021148036253776
$ws.Range("L12").Value = 0.03021148036253777
$ws.Range("S12").Value = 0.03625377643504532
$ws.Range("G13").Value = 0.7169811320754716
$ws.Range("J13").Value = 0.2264150943396226
$ws.Range("S13").Value = 0.05660377358490566
$ws.Range("F15").Value = 0.01762114537444934
$ws.Range("H15").Value = 0.1585903083700441
$ws.Range("I15").Value = 0.0881057268722467
$ws.Range("J15").Value = 0.3744493392070485
$ws.Range("K15").Value = 0.03524229074889868
$ws.Range("M15").Value = 0.013215859030837
$ws.Range("N15").Value = 0.002202643171806168
$ws.Range("O15").Value = 0.05947136563876652
$ws.Range("S15").Value = 0.2511013215859031
$ws.Range("F16").Value = 0.0178117048346056
$ws.Range("H16").Value = 0.1552162849872774
$ws.Range("I16").Value = 0.06870229007633588
$ws.Range("J16").Value = 0.4478371501272265
$ws.Range("K16").Value = 0.09923664122137404
$ws.Range("M16").Value = 0.02290076335877863
$ws.Range("N16").Value = 0.005089058524173028
$ws.Range("O16").Value = 0.04580152671755725
$ws.Range("S16").Value = 0.1374045801526718
$ws.Range("F17").Value = 0.01203501094091904
$ws.Range("H17").Value = 0.1892778993435449
$ws.Range("I17").Value = 0.1072210065645514
$ws.Range("J17").Value = 0.3982494529540481
$ws.Range("K17").Value = 0.0962800875273523
$ws.Range("M17").Value = 0.01531728665207877
$ws.Range("N17").Value = 0.001094091903719912
$ws.Range("O17").Value = 0.06236323851203501
$ws.Range("S17").Value = 0.1181619256017506
$ws.Range("F18").Value = 0.0128755364806867
$ws.Range("H18").Value = 0.1802575107296137
$ws.Range("I18").Value = 0.09012875536480687
$ws.Range("J18").Value = 0.4248927038626609
$ws.Range("K18").Value = 0.08798283261802575
$ws.Range("M18").Value = 0.01931330472103004
$ws.Range("O18").Value = 0.06866952789699571
$ws.Range("S18").Value = 0.1158798283261803
$ws.Range("F19").Value = 0.01658433309809457
$ws.Range("H19").Value = 0.2275935074100212
$ws.Range("I19").Value = 0.08997882851093861
$ws.Range("J19").Value = 0.3599153140437544
$ws.Range("K19").Value = 0.09562455892731123
$ws.Range("M19").Value = 0.02540578687367678
$ws.Range("N19").Value = 0.0003528581510232887
$ws.Range("O19").Value = 0.06563161609033169
$ws.Range("S19").Value = 0.1189131968948483
